$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the new "Priority" header in D1
$ws.Range("D1").Value = "Priority"

# Priority values for rows 2..91 (column D), grouped in blocks of 9 rows
# matching the route groups already present in column A.
$priorities = @(0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,2,2,2,2,2,2,2,2,2,1,1,1,1,1,1,1,1,1,1,1,1,1,1,1,1,1,1,1,1,1,1,1,1,1,1,1,3,3,3,3,3,3,3,3,3,2,2,2,2,2,2,2,2,2,3,3,3,3,3,3,3,3,3,2,2,2,2,2,2,2,2,2)

for ($i = 0; $i -lt $priorities.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 4).Value = $priorities[$i]
}

# Update the selection / view state to mirror the saved workbook view
$ws.Range("F22").Select()
$excel.ActiveWindow.ScrollRow = 8
$excel.ActiveWindow.ScrollColumn = 1
